$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number stamped on every data
# row (2 through 452). This refresh bumps that serial by one day,
# from 45189 (2023-09-20) to 45190 (2023-09-21), for every row.
$range = $ws.Range("C2:C452")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
